$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 11:50"

$data = @(
  @("Estados Unidos", 188592, 62, 7251, 177285, 4576, 3, 4056),
  @("Italia", 105792, 0, 15729, 77635, 4023, 0, 12428),
  @("España", 102136, 6213, 22647, 70436, 5872, 589, 9053),
  @("China", 81554, 36, 76238, 2004, 466, 7, 3312),
  @("Alemania", 71921, 113, 16100, 55040, 2675, 6, 781),
  @("Francia", 52128, 0, 9444, 39161, 5565, 0, 3523),
  @("Iran", 44605, 0, 14656, 27051, 3703, 0, 2898),
  @("Reino Unido", 25150, 0, 135, 23226, 163, 0, 1789),
  @("Suiza", 16605, 0, 1823, 14349, 301, 0, 433),
  @("Belgica", 13964, 1189, 2132, 11004, 1088, 123, 828),
  @("Turquia", 13531, 0, 243, 13074, 847, 0, 214),
  @("Paises Bajos", 12595, 0, 250, 11306, 1053, 0, 1039),
  @("Austria", 10366, 186, 1095, 9125, 215, 18, 146),
  @("Corea del Sur", 9887, 101, 5567, 4155, 55, 3, 165),
  @("Canada", 8612, 0, 1242, 7269, 120, 0, 101),
  @("Portugal", 7443, 0, 43, 7240, 188, 0, 160),
  @("Brasil", 5812, 95, 127, 5483, 296, 1, 202),
  @("Israel", 5591, 233, 226, 5344, 97, 1, 21),
  @("Australia", 4862, 99, 345, 4496, 50, 1, 21),
  @("Noruega", 4651, 10, 13, 4599, 97, 0, 39),
  @("Suecia", 4435, 0, 16, 4239, 358, 0, 180),
  @("Chequia", 3330, 22, 45, 3253, 70, 1, 32),
  @("Irlanda", 3235, 0, 5, 3159, 103, 0, 71),
  @("Malasia", 2908, 142, 645, 2218, 102, 2, 45),
  @("Dinamarca", 2860, 0, 1, 2769, 145, 0, 90),
  @("Rusia", 2777, 440, 190, 2563, 8, 7, 24),
  @("Chile", 2738, 0, 156, 2570, 14, 0, 12),
  @("Polonia", 2347, 36, 7, 2305, 50, 2, 35),
  @("Filipinas", 2311, 227, 50, 2165, 1, 8, 96),
  @("Ecuador", 2302, 0, 58, 2165, 100, 0, 79),
  @("Rumania", 2245, 0, 220, 1940, 62, 3, 85),
  @("Luxemburgo", 2178, 0, 80, 2075, 31, 0, 23),
  @("Japon", 2178, 0, 424, 1697, 69, 0, 57),
  @("Pakistan", 2042, 104, 82, 1934, 12, 0, 26),
  @("Tailandia", 1771, 120, 416, 1343, 23, 2, 12),
  @("Indonesia", 1677, 149, 103, 1417, 0, 21, 157),
  @("India", 1590, 193, 148, 1397, 0, 10, 45),
  @("Arabia Saudita", 1563, 0, 165, 1388, 31, 0, 10),
  @("Finlandia", 1418, 0, 10, 1391, 56, 0, 17),
  @("Sudafrica", 1353, 0, 50, 1298, 7, 0, 5),
  @("Grecia", 1314, 0, 52, 1213, 72, 0, 49),
  @("Mexico", 1215, 121, 35, 1151, 1, 1, 29),
  @("Panama", 1181, 0, 9, 1142, 50, 0, 30),
  @("Islandia", 1135, 0, 198, 935, 11, 0, 2),
  @("Republica Dominicana", 1109, 0, 5, 1053, 0, 0, 51),
  @("Peru", 1065, 0, 394, 641, 49, 0, 30),
  @("Argentina", 1054, 0, 240, 787, 0, 0, 27),
  @("Singapur", 926, 0, 240, 683, 22, 0, 3),
  @("Colombia", 906, 0, 31, 859, 35, 0, 16),
  @("Serbia", 900, 0, 42, 835, 62, 0, 23),
  @("Croacia", 867, 0, 67, 794, 32, 0, 6),
  @("Eslovenia", 841, 39, 10, 816, 31, 0, 15),
  @("Catar", 781, 0, 62, 717, 6, 0, 2),
  @("Estonia", 779, 34, 26, 749, 13, 0, 4),
  @("Hong Kong", 765, 50, 147, 614, 5, 0, 4),
  @("Argelia", 716, 0, 46, 626, 0, 0, 44),
  @("Crucero", 712, 0, 603, 98, 15, 0, 11),
  @("Egipto", 710, 0, 157, 507, 0, 0, 46),
  @("Nueva Zelanda", 708, 61, 83, 624, 2, 0, 1),
  @("Irak", 694, 0, 170, 474, 0, 0, 50),
  @("Ucrania", 669, 24, 10, 642, 0, 0, 17),
  @("Emiratos Arabes Unidos", 664, 0, 61, 597, 2, 0, 6),
  @("Marruecos", 638, 21, 24, 578, 1, 0, 36),
  @("Lituania", 581, 44, 7, 566, 27, 0, 8),
  @("Armenia", 571, 39, 31, 537, 30, 0, 3),
  @("Barein", 567, 0, 316, 247, 2, 0, 4),
  @("Hungria", 525, 33, 40, 465, 6, 4, 20),
  @("Libano", 479, 16, 37, 430, 7, 0, 12),
  @("Bosnia y Herzegovina", 448, 28, 17, 418, 1, 0, 13),
  @("Letonia", 446, 48, 1, 445, 3, 0, 0),
  @("Bulgaria", 412, 13, 20, 384, 14, 0, 8),
  @("Eslovaquia", 400, 37, 3, 397, 1, 0, 0),
  @("Tunez", 394, 0, 3, 381, 10, 0, 10),
  @("Principado de Andorra", 376, 0, 10, 354, 12, 0, 12),
  @("Kazajistan", 369, 26, 24, 342, 6, 1, 3),
  @("Moldavia", 353, 0, 22, 327, 44, 0, 4),
  @("Costa Rica", 347, 0, 4, 341, 8, 0, 2),
  @("Uruguay", 338, 0, 41, 296, 12, 0, 1),
  @("Republica de Macedonia", 329, 0, 12, 308, 4, 0, 9),
  @("Taiwan", 329, 7, 45, 279, 0, 0, 5),
  @("Kuwait", 317, 28, 80, 237, 13, 0, 0),
  @("Azerbaiyan", 298, 0, 26, 267, 11, 0, 5),
  @("Jordania", 274, 0, 30, 239, 5, 0, 5),
  @("Republica de Chipre", 262, 0, 23, 231, 3, 0, 8),
  @("Burkina Faso", 261, 0, 32, 215, 0, 0, 14),
  @("Reunion", 247, 0, 1, 246, 4, 0, 0),
  @("Albania", 243, 0, 52, 176, 8, 0, 15),
  @("San Marino", 236, 0, 13, 197, 16, 0, 26),
  @("Camerun", 223, 30, 5, 212, 0, 0, 6),
  @("Vietnam", 212, 0, 58, 154, 3, 0, 0),
  @("Oman", 210, 18, 34, 175, 3, 0, 1),
  @("Afganistan", 196, 22, 5, 187, 0, 0, 4),
  @("Cuba", 186, 0, 8, 172, 3, 0, 6),
  @("Costa de Marfil", 179, 0, 7, 171, 0, 0, 1),
  @("Senegal", 175, 0, 40, 134, 0, 1, 1),
  @("Uzbekistan", 173, 1, 8, 163, 8, 0, 2),
  @("Islas Feroe", 173, 4, 75, 98, 1, 0, 0),
  @("Honduras", 172, 31, 3, 159, 4, 3, 10),
  @("Malta", 169, 0, 2, 167, 2, 0, 0),
  @("Bielorrusia", 163, 11, 53, 109, 2, 0, 1),
  @("Ghana", 161, 0, 31, 125, 1, 0, 5),
  @("Mauricio", 147, 4, 0, 142, 1, 0, 5),
  @("Sri Lanka", 143, 0, 18, 123, 5, 0, 2),
  @("Venezuela", 143, 8, 41, 99, 6, 0, 3),
  @("Nigeria", 139, 4, 9, 128, 0, 0, 2),
  @("Estado de Palestina", 134, 15, 18, 115, 0, 0, 1),
  @("Brunei", 129, 0, 45, 83, 3, 0, 1),
  @("Martinica", 128, 0, 27, 98, 15, 0, 3),
  @("Montenegro", 120, 11, 0, 118, 4, 0, 2),
  @("Bolivia", 115, 8, 1, 107, 3, 1, 7),
  @("Georgia", 115, 5, 23, 92, 6, 0, 0),
  @("Guadalupe", 114, 0, 22, 88, 14, 0, 4),
  @("Kirguistan", 111, 4, 3, 108, 3, 0, 0),
  @("Consejo Danes para los Refugiados", 109, 11, 2, 99, 0, 0, 8),
  @("Camboya", 109, 0, 25, 84, 1, 0, 0),
  @("Mayotte", 94, 0, 10, 83, 3, 0, 1),
  @("Trinidad yTobago", 89, 2, 1, 84, 0, 1, 4),
  @("Ruanda", 75, 0, 0, 75, 0, 0, 0),
  @("Paraguay", 69, 4, 1, 65, 3, 0, 3),
  @("Gibraltar", 69, 0, 34, 35, 0, 0, 0),
  @("Liechtenstein", 68, 0, 0, 68, 0, 0, 0),
  @("Isla de Man", 65, 5, 0, 65, 0, 0, 0),
  @("Kenia", 59, 0, 1, 57, 2, 0, 1),
  @("Madagascar", 57, 0, 0, 57, 6, 0, 0),
  @("Aruba", 55, 0, 1, 54, 0, 0, 0),
  @("Banglades", 54, 3, 25, 23, 1, 1, 6),
  @("Monaco", 52, 0, 2, 49, 0, 0, 1),
  @("Guayana Francesa", 51, 0, 15, 36, 0, 0, 0),
  @("Uganda", 44, 0, 0, 44, 0, 0, 0),
  @("Macao", 41, 0, 10, 31, 0, 0, 0),
  @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2),
  @("Guatemala", 39, 1, 12, 26, 1, 0, 1),
  @("Jamaica", 38, 0, 2, 34, 0, 0, 2),
  @("Polinesia Francesa", 37, 0, 0, 37, 1, 0, 0),
  @("Zambia", 36, 0, 0, 36, 0, 0, 0),
  @("Barbados", 34, 0, 0, 34, 0, 0, 0),
  @("Niger", 34, 0, 0, 31, 0, 0, 3),
  @("Togo", 34, 0, 10, 23, 0, 0, 1),
  @("Guam", 32, 0, 0, 31, 0, 0, 1),
  @("El Salvador", 32, 0, 0, 31, 1, 0, 1),
  @("Bermudas", 32, 0, 10, 22, 0, 0, 0),
  @("Republica de Yibuti", 30, 0, 0, 30, 0, 0, 0),
  @("Etiopia", 29, 3, 2, 27, 2, 0, 0),
  @("Mali", 28, 0, 0, 26, 0, 0, 2),
  @("Guinea", 22, 0, 0, 22, 0, 0, 0),
  @("Congo", 19, 0, 0, 19, 0, 0, 0),
  @("Tanzania", 19, 0, 1, 17, 0, 0, 1),
  @("Gabon", 18, 2, 0, 17, 0, 0, 1),
  @("Maldivas", 18, 0, 13, 5, 0, 0, 0),
  @("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0),
  @("Nueva Caledonia", 16, 0, 0, 16, 0, 0, 0),
  @("Haiti", 16, 1, 1, 15, 0, 0, 0),
  @("San Martin (Parte Francesa)", 16, 1, 2, 13, 0, 0, 1),
  @("Eritrea", 15, 0, 0, 15, 0, 0, 0),
  @("Guinea Ecuatorial", 15, 0, 1, 14, 0, 0, 0),
  @("Birmania", 15, 0, 0, 14, 0, 0, 1),
  @("Bahamas", 15, 1, 1, 14, 0, 0, 0),
  @("Islas Caimanes", 14, 0, 0, 13, 0, 0, 1),
  @("Mongolia", 14, 2, 2, 12, 0, 0, 0),
  @("Santa Lucia", 13, 0, 1, 12, 0, 0, 0),
  @("Dominica", 12, 0, 0, 12, 0, 0, 0),
  @("Guyana", 12, 0, 0, 10, 0, 0, 2),
  @("Namibia", 11, 0, 2, 9, 0, 0, 0),
  @("Curazao", 11, 0, 2, 8, 0, 0, 1),
  @("Seychelles", 10, 0, 0, 10, 0, 0, 0),
  @("Laos", 10, 1, 0, 10, 0, 0, 0),
  @("Surinam", 10, 1, 0, 10, 0, 0, 0),
  @("Libia", 10, 0, 1, 9, 0, 0, 0),
  @("Groenlandia", 10, 0, 2, 8, 0, 0, 0),
  @("Siria", 10, 0, 0, 8, 0, 0, 2),
  @("Suazilandia", 9, 0, 0, 9, 0, 0, 0),
  @("Granada", 9, 0, 0, 9, 0, 0, 0),
  @("Benin", 9, 0, 1, 8, 0, 0, 0),
  @("San Cristobal y Nieves", 8, 0, 0, 8, 0, 0, 0),
  @("Guinea-Bisau", 8, 0, 0, 8, 0, 0, 0),
  @("Mozambique", 8, 0, 0, 8, 0, 0, 0),
  @("Zimbabue", 8, 0, 0, 7, 0, 0, 1),
  @("Antigua y Barbuda", 7, 0, 0, 7, 0, 0, 0),
  @("Republica del Chad", 7, 0, 0, 7, 0, 0, 0),
  @("Angola", 7, 0, 1, 4, 0, 0, 2),
  @("Sudan", 7, 0, 1, 4, 0, 0, 2),
  @("San Martin (Parte Holandesa)", 6, 0, 0, 6, 0, 0, 0),
  @("Liberia", 6, 3, 0, 6, 0, 0, 0),
  @("Santa Sede", 6, 0, 0, 6, 0, 0, 0),
  @("Cabo Verde", 6, 0, 0, 5, 0, 0, 1),
  @("San Bartolome", 6, 0, 1, 5, 0, 0, 0),
  @("Mauritania", 6, 0, 2, 3, 0, 0, 1),
  @("Islas Turcas y Caicos", 5, 0, 0, 5, 0, 0, 0),
  @("Fiyi", 5, 0, 0, 5, 0, 0, 0),
  @("Montserrat", 5, 0, 0, 5, 0, 0, 0),
  @("Somalia", 5, 0, 1, 4, 0, 0, 0),
  @("Nepal", 5, 0, 1, 4, 0, 0, 0),
  @("Nicaragua", 5, 0, 0, 4, 0, 0, 1),
  @("Butan", 4, 0, 0, 4, 0, 0, 0),
  @("Gambia", 4, 0, 0, 3, 0, 0, 1),
  @("Botsuana", 4, 0, 0, 3, 0, 0, 1),
  @("Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0),
  @("Belice", 3, 0, 0, 3, 0, 0, 0),
  @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
  @("Bonaire, San Eustaquio y Saba", 2, 2, 0, 2, 0, 0, 0),
  @("Anguila", 2, 0, 0, 2, 0, 0, 0),
  @("Burundi", 2, 0, 0, 2, 0, 0, 0),
  @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
  @("Sierra Leona", 1, 0, 0, 1, 0, 0, 0),
  @("Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
  @("San Vicente y las Granadinas", 1, 0, 1, 0, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

"Done: " + $data.Length
